$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new data values: C2 (Planned Velocity, Iteration 2) = 40,
# B3 (Actual Velocity, Iteration 1) = 35
$ws.Range("C2").Value = 40
$ws.Range("B3").Value = 35

# Update the active selection to E11 (matches saved selection state in diff)
$ws.Range("E11").Select()
